{"js": "// Word JS API (Office.js) edit script.\n// Body is the implementation of: async (context) => { ... }\n\nconst replacements = [\n  [\n    \"Ativa\u00e7\u00e3o: 01/01/2021\",\n    \"Ativa\u00e7\u00e3o: 01/01/2024\",\n  ],\n  [\n    \"Apresentar no\u00e7\u00f5es de Matem\u00e1tica Financeira, Gest\u00e3o Financeiras e Engenharia Econ\u00f4mica aos alunos, capacitando-os para construir e analisar fluxos de caixa de projetos e empreendimentos, discutir os principais aspectos da gest\u00e3o financeira das empresas industriais, comerciais e de servi\u00e7os e analisar e propor estrat\u00e9gias de gest\u00e3o financeira relacionadas \u00e0s estrat\u00e9gias de mercado e de produ\u00e7\u00e3o.\",\n    \"Apresentar conceitos b\u00e1sicos de Engenharia Econ\u00f4mica e Gest\u00e3o Financeira aos alunos, capacitando-os para analisar projetos de investimentos e avaliar o desempenho financeiro de organiza\u00e7\u00f5es.\",\n  ],\n  [\n    \"To present notions of Financial Mathematics, Financial Management and Economic Engineering to students, enabling them to build and analyze cash flows from projects and enterprises, discuss the main aspects of financial management in industrial, commercial and service companies and analyze and propose strategies for financial management related to market and production strategies.\",\n    \"Introduce basic concepts of Economic Engineering and Financial Management to students, preparing them to analyze investment projects and evaluate the financial performance of organizations.\",\n  ],\n  [\n    \"1. Engenharia Econ\u00f4mica. 2 \u2013 Finan\u00e7as\",\n    \"A. Engenharia Econ\u00f4mica. B. Finan\u00e7as.\",\n  ],\n  [\n    \"1. Economic Engineering. 2 \u2013 Finance\",\n    \"A. Economic Engineering. B. Finance.\",\n  ],\n  [\n    \"1. Engenharia Econ\u00f4mica: Vari\u00e1vel tempo: juros simples, juros compostos; M\u00e9todos de amortiza\u00e7\u00e3o; Equival\u00eancia de m\u00e9todos; M\u00e9todos de Decis\u00e3o; Renova\u00e7\u00e3o e substitui\u00e7\u00e3o de equipamentos; Deprecia\u00e7\u00e3o; An\u00e1lise de Projetos, Riscos em projetos; Estimativa do custo de capital pr\u00f3prio (CAPM) e WACC.2. Finan\u00e7as: O ciclo da produ\u00e7\u00e3o e o ciclo do capital; An\u00e1lise de \u00cdndices; Fontes de Financiamento, Alavancagem; Capital de Giro; Custo de Capital; A\u00e7\u00f5es, Pol\u00edtica de Dividendos; Financiamento de Longo Prazo, Corporate Finance/Project Finance; EVA e MVA.\",\n    \"A. Engenharia Econ\u00f4mica: s\u00e9rie de pagamentos; juros simples; juros compostos; equival\u00eancia de taxas; m\u00e9todos de amortiza\u00e7\u00e3o e planos de pagamentos; m\u00e9todos de avalia\u00e7\u00e3o de projetos de investimentos; introdu\u00e7\u00e3o a avalia\u00e7\u00e3o de riscos. B. Finan\u00e7as: o ciclo da produ\u00e7\u00e3o e o ciclo do capital; financiamento; alavancagem; capital de giro; custo de capital; pol\u00edticas de dividendos; Indicadores de desempenho financeiro.\",\n  ],\n  [\n    \"1. Economic Engineering: Time variable: simple interest, compound interest; Amortization methods; Equivalence of methods; Decision Methods; Renovation and replacement of equipment; Depreciation; Project Analysis, Project Risks; Estimated cost of equity (CAPM) and WACC.2. Finance: The production cycle and the capital cycle; Index Analysis; Financing Sources, Leverage; Working capital; Capital cost; Shares, Dividend Policy; Long Term Financing, Corporate Finance/Project Finance; EVA and MVA.\",\n    \"A. Economic Engineering: series of payments; simple interest; compound interest; fee equivalence; amortization methods and payment schedules; investment project evaluation methods; introduction to risk assessment. B. Finance: the production cycle and the capital cycle; financing; leverage; working capital; cost of capital; dividend policies; financial performance indicators.\",\n  ],\n  [\n    \"Aulas expositivas e dialogadas; din\u00e2micas, projetos e trabalhos em grupo; exerc\u00edcios individuais; e, semin\u00e1rios, debates e palestras.\",\n    \"Provas, trabalhos em grupo, exerc\u00edcios individuais, e semin\u00e1rios.\",\n  ],\n  [\n    \"M\u00e9dia Aritm\u00e9tica dos Projetos, Trabalhos, Exerc\u00edcios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as quest\u00f5es relativas \u00e0s Compet\u00eancias (Conhecimento, Habilidade e Atitude, que incluem a presen\u00e7a e participa\u00e7\u00e3o dos alunos nas aulas) desenvolvidas\",\n    \"M\u00e9dia das atividades avaliativas.\",\n  ],\n  [\n    \"GITMAN, L. J. - ZUTTER, C. J. Princ\u00edpios de Administra\u00e7\u00e3o Financeira. 14 ed. S\u00e3o Paulo: Perason, 2017.GROPPELLI, A. A.; NIKBAKHT, E. Administra\u00e7\u00e3o Financeira. 3 ed. S\u00e3o Paulo: Saraiva, 2010.MARCOUS\u00c9, I.; SURRIDGE, M.; GILLESPIE, A. Finan\u00e7as. S\u00e3o Paulo: Saraiva, 2013.ASSAF NETO, A. E LIMA, F. G. 3 ed. CURSO DE ADMINISTRA\u00c7\u00c3O FINANCEIRA. S\u00e3o Paulo: Atlas, 2014MARIANO, F.; MENESES, A. Curso De Administra\u00e7\u00e3o Financeira. S\u00e3o Paulo: M\u00e9todo, 2012.MORANTE, A. S. An\u00e1lise das Demonstra\u00e7\u00f5es Financeiras. 2 ed. S\u00e3o Paulo: Atlas, 2009.NEWNAN, D. G.; LAVELLE, J. P. Fundamentos de Engenharia Econ\u00f4mica. S\u00e3o Paulo: LTC,2000.KOPITTKE, B. H.; CASAROTTO FILHO, N. AN\u00c1LISE DE INVESTIMENTOS: Matem\u00e1tica Financeira, Engenharia Econ\u00f4mica, Estrat\u00e9gia Empresarial. 11 ed. S\u00e3o Paulo: Atlas, 2010.HOJI, M.; LUZ, A. E. Gest\u00e3o Financeira Econ\u00f4mica: Did\u00e1tica, Objetiva e Pr\u00e1tica. S\u00e3o Paulo: Atlas, 2019.GOMES, J. M. Elabora\u00e7\u00e3o e An\u00e1lise De Viabilidade Econ\u00f4mica De Projetos. S\u00e3o Paulo: Atlas, 2013.OLIVO, R. L. F. An\u00e1lise de Investimentos. Campinas: Al\u00ednea, 2011.ALMEIDA, J. T. S. Matem\u00e1tica Financeira. Rio de Janeiro: LTC, 2016.EHRLICH, Pierre Jacques. Engenharia Econ\u00f4mica. S\u00e3o Paulo: Editora Atlas, 2005.HIRSCHFELD, Henrique. Engenharia econ\u00f4mica e an\u00e1lise de custos. 7. ed. S\u00e3o Paulo: Atlas, 2007.MOTTA, Regis da Rocha; CAL\u00d4BA, Guilherme Marques. An\u00e1lise de Investimentos. S\u00e3o Paulo: Atlas 2002.SANVICENTE, A. Z. Administra\u00e7\u00e3o Financeira. S\u00e3o Paulo: Editora Atlas, 2007.VAN HORNE, J. C. Pol\u00edtica e Administra\u00e7\u00e3o Financeira. Rio de Janeiro: Livros T\u00e9cnicos e Cient\u00edficos, 1974.WESTON, J. F.; BRIGHAM, E. F. Administra\u00e7\u00e3o Financeira de Empresas. S\u00e3o Paulo: Editora Interamericana, 2000.\",\n    \"EHRLICH, P. J. Engenharia Econ\u00f4mica. S\u00e3o Paulo: Editora Atlas, 2005.JAFFE, R. W. Administra\u00e7\u00e3o Financeira. 2 ed. S\u00e3o Paulo, Editora Atlas, 2002.GITMAN, L. J. - ZUTTER, C. J. Princ\u00edpios de Administra\u00e7\u00e3o Financeira. 14 ed. S\u00e3o Paulo: Pearson, 2017.ASSAF NETO, A. E LIMA, F. G. 3 ed. CURSO DE ADMINISTRA\u00c7\u00c3O FINANCEIRA. S\u00e3o Paulo: Atlas, 2014MORANTE, A. S. An\u00e1lise das Demonstra\u00e7\u00f5es Financeiras. 2 ed. S\u00e3o Paulo: Atlas, 2009.NEWNAN, D. G.; LAVELLE, J. P. Fundamentos de Engenharia Econ\u00f4mica. S\u00e3o Paulo: LTC, 2000.HOJI, M.; LUZ, A. E. Gest\u00e3o Financeira Econ\u00f4mica: Did\u00e1tica, Objetiva e Pr\u00e1tica. S\u00e3o Paulo: Atlas, 2019.SANVICENTE, A. Z. Administra\u00e7\u00e3o Financeira. S\u00e3o Paulo: Editora Atlas, 2007.\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText.substring(0, 60));\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop edit script.\n# $word / $app -> Word.Application ; $d -> ActiveDocument (also exposed as $word.ActiveDocument)\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $found = $range.Find.Execute(\n        $findText,      # FindText\n        $false,         # MatchCase\n        $true,          # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap (wdFindContinue)\n        $false,         # Format\n        $replaceText,   # ReplaceWith\n        2               # Replace (wdReplaceOne)\n    )\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\nReplace-ExactText \"Ativa\u00e7\u00e3o: 01/01/2021\" \"Ativa\u00e7\u00e3o: 01/01/2024\"\n\nReplace-ExactText \"Apresentar no\u00e7\u00f5es de Matem\u00e1tica Financeira, Gest\u00e3o Financeiras e Engenharia Econ\u00f4mica aos alunos, capacitando-os para construir e analisar fluxos de caixa de projetos e empreendimentos, discutir os principais aspectos da gest\u00e3o financeira das empresas industriais, comerciais e de servi\u00e7os e analisar e propor estrat\u00e9gias de gest\u00e3o financeira relacionadas \u00e0s estrat\u00e9gias de mercado e de produ\u00e7\u00e3o.\" \"Apresentar conceitos b\u00e1sicos de Engenharia Econ\u00f4mica e Gest\u00e3o Financeira aos alunos, capacitando-os para analisar projetos de investimentos e avaliar o desempenho financeiro de organiza\u00e7\u00f5es.\"\n\nReplace-ExactText \"To present notions of Financial Mathematics, Financial Management and Economic Engineering to students, enabling them to build and analyze cash flows from projects and enterprises, discuss the main aspects of financial management in industrial, commercial and service companies and analyze and propose strategies for financial management related to market and production strategies.\" \"Introduce basic concepts of Economic Engineering and Financial Management to students, preparing them to analyze investment projects and evaluate the financial performance of organizations.\"\n\nReplace-ExactText \"1. Engenharia Econ\u00f4mica. 2 \u2013 Finan\u00e7as\" \"A. Engenharia Econ\u00f4mica. B. Finan\u00e7as.\"\n\nReplace-ExactText \"1. Economic Engineering. 2 \u2013 Finance\" \"A. Economic Engineering. B. Finance.\"\n\nReplace-ExactText \"1. Engenharia Econ\u00f4mica: Vari\u00e1vel tempo: juros simples, juros compostos; M\u00e9todos de amortiza\u00e7\u00e3o; Equival\u00eancia de m\u00e9todos; M\u00e9todos de Decis\u00e3o; Renova\u00e7\u00e3o e substitui\u00e7\u00e3o de equipamentos; Deprecia\u00e7\u00e3o; An\u00e1lise de Projetos, Riscos em projetos; Estimativa do custo de capital pr\u00f3prio (CAPM) e WACC.2. Finan\u00e7as: O ciclo da produ\u00e7\u00e3o e o ciclo do capital; An\u00e1lise de \u00cdndices; Fontes de Financiamento, Alavancagem; Capital de Giro; Custo de Capital; A\u00e7\u00f5es, Pol\u00edtica de Dividendos; Financiamento de Longo Prazo, Corporate Finance/Project Finance; EVA e MVA.\" \"A. Engenharia Econ\u00f4mica: s\u00e9rie de pagamentos; juros simples; juros compostos; equival\u00eancia de taxas; m\u00e9todos de amortiza\u00e7\u00e3o e planos de pagamentos; m\u00e9todos de avalia\u00e7\u00e3o de projetos de investimentos; introdu\u00e7\u00e3o a avalia\u00e7\u00e3o de riscos. B. Finan\u00e7as: o ciclo da produ\u00e7\u00e3o e o ciclo do capital; financiamento; alavancagem; capital de giro; custo de capital; pol\u00edticas de dividendos; Indicadores de desempenho financeiro.\"\n\nReplace-ExactText \"1. Economic Engineering: Time variable: simple interest, compound interest; Amortization methods; Equivalence of methods; Decision Methods; Renovation and replacement of equipment; Depreciation; Project Analysis, Project Risks; Estimated cost of equity (CAPM) and WACC.2. Finance: The production cycle and the capital cycle; Index Analysis; Financing Sources, Leverage; Working capital; Capital cost; Shares, Dividend Policy; Long Term Financing, Corporate Finance/Project Finance; EVA and MVA.\" \"A. Economic Engineering: series of payments; simple interest; compound interest; fee equivalence; amortization methods and payment schedules; investment project evaluation methods; introduction to risk assessment. B. Finance: the production cycle and the capital cycle; financing; leverage; working capital; cost of capital; dividend policies; financial performance indicators.\"\n\nReplace-ExactText \"Aulas expositivas e dialogadas; din\u00e2micas, projetos e trabalhos em grupo; exerc\u00edcios individuais; e, semin\u00e1rios, debates e palestras.\" \"Provas, trabalhos em grupo, exerc\u00edcios individuais, e semin\u00e1rios.\"\n\nReplace-ExactText \"M\u00e9dia Aritm\u00e9tica dos Projetos, Trabalhos, Exerc\u00edcios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as quest\u00f5es relativas \u00e0s Compet\u00eancias (Conhecimento, Habilidade e Atitude, que incluem a presen\u00e7a e participa\u00e7\u00e3o dos alunos nas aulas) desenvolvidas\" \"M\u00e9dia das atividades avaliativas.\"\n\nReplace-ExactText \"GITMAN, L. J. - ZUTTER, C. J. Princ\u00edpios de Administra\u00e7\u00e3o Financeira. 14 ed. S\u00e3o Paulo: Perason, 2017.GROPPELLI, A. A.; NIKBAKHT, E. Administra\u00e7\u00e3o Financeira. 3 ed. S\u00e3o Paulo: Saraiva, 2010.MARCOUS\u00c9, I.; SURRIDGE, M.; GILLESPIE, A. Finan\u00e7as. S\u00e3o Paulo: Saraiva, 2013.ASSAF NETO, A. E LIMA, F. G. 3 ed. CURSO DE ADMINISTRA\u00c7\u00c3O FINANCEIRA. S\u00e3o Paulo: Atlas, 2014MARIANO, F.; MENESES, A. Curso De Administra\u00e7\u00e3o Financeira. S\u00e3o Paulo: M\u00e9todo, 2012.MORANTE, A. S. An\u00e1lise das Demonstra\u00e7\u00f5es Financeiras. 2 ed. S\u00e3o Paulo: Atlas, 2009.NEWNAN, D. G.; LAVELLE, J. P. Fundamentos de Engenharia Econ\u00f4mica. S\u00e3o Paulo: LTC,2000.KOPITTKE, B. H.; CASAROTTO FILHO, N. AN\u00c1LISE DE INVESTIMENTOS: Matem\u00e1tica Financeira, Engenharia Econ\u00f4mica, Estrat\u00e9gia Empresarial. 11 ed. S\u00e3o Paulo: Atlas, 2010.HOJI, M.; LUZ, A. E. Gest\u00e3o Financeira Econ\u00f4mica: Did\u00e1tica, Objetiva e Pr\u00e1tica. S\u00e3o Paulo: Atlas, 2019.GOMES, J. M. Elabora\u00e7\u00e3o e An\u00e1lise De Viabilidade Econ\u00f4mica De Projetos. S\u00e3o Paulo: Atlas, 2013.OLIVO, R. L. F. An\u00e1lise de Investimentos. Campinas: Al\u00ednea, 2011.ALMEIDA, J. T. S. Matem\u00e1tica Financeira. Rio de Janeiro: LTC, 2016.EHRLICH, Pierre Jacques. Engenharia Econ\u00f4mica. S\u00e3o Paulo: Editora Atlas, 2005.HIRSCHFELD, Henrique. Engenharia econ\u00f4mica e an\u00e1lise de custos. 7. ed. S\u00e3o Paulo: Atlas, 2007.MOTTA, Regis da Rocha; CAL\u00d4BA, Guilherme Marques. An\u00e1lise de Investimentos. S\u00e3o Paulo: Atlas 2002.SANVICENTE, A. Z. Administra\u00e7\u00e3o Financeira. S\u00e3o Paulo: Editora Atlas, 2007.VAN HORNE, J. C. Pol\u00edtica e Administra\u00e7\u00e3o Financeira. Rio de Janeiro: Livros T\u00e9cnicos e Cient\u00edficos, 1974.WESTON, J. F.; BRIGHAM, E. F. Administra\u00e7\u00e3o Financeira de Empresas. S\u00e3o Paulo: Editora Interamericana, 2000.\" \"EHRLICH, P. J. Engenharia Econ\u00f4mica. S\u00e3o Paulo: Editora Atlas, 2005.JAFFE, R. W. Administra\u00e7\u00e3o Financeira. 2 ed. S\u00e3o Paulo, Editora Atlas, 2002.GITMAN, L. J. - ZUTTER, C. J. Princ\u00edpios de Administra\u00e7\u00e3o Financeira. 14 ed. S\u00e3o Paulo: Pearson, 2017.ASSAF NETO, A. E LIMA, F. G. 3 ed. CURSO DE ADMINISTRA\u00c7\u00c3O FINANCEIRA. S\u00e3o Paulo: Atlas, 2014MORANTE, A. S. An\u00e1lise das Demonstra\u00e7\u00f5es Financeiras. 2 ed. S\u00e3o Paulo: Atlas, 2009.NEWNAN, D. G.; LAVELLE, J. P. Fundamentos de Engenharia Econ\u00f4mica. S\u00e3o Paulo: LTC, 2000.HOJI, M.; LUZ, A. E. Gest\u00e3o Financeira Econ\u00f4mica: Did\u00e1tica, Objetiva e Pr\u00e1tica. S\u00e3o Paulo: Atlas, 2019.SANVICENTE, A. Z. Administra\u00e7\u00e3o Financeira. S\u00e3o Paulo: Editora Atlas, 2007.\"\n"}
